# İş Takip Güncellemesi - 23.02.2026 14:36:51

$wb = $excel.ActiveWorkbook

# --- 1) "İş Takip Listesi" sheet: update status for row 117 ---
$wsTakip = $wb.Worksheets.Item("İş Takip Listesi")
$wsTakip.Range("L117").Value = "BİLGİLENDİRME İLANINDA"

# --- 2) "Güncelleme" sheet: update evaluation status + add bilgilendirme ilanı date for row 24 ---
$wsGuncelleme = $wb.Worksheets.Item("Güncelleme")
$wsGuncelleme.Range("K24").Value = "Yapıldı"

# M24 holds a date formatted as plain text ("YYYY-MM-DD"), like the rest of column M.
# Force text formatting first so Excel doesn't reinterpret the string as a date serial.
$cellM24 = $wsGuncelleme.Range("M24")
$cellM24.NumberFormat = "@"
$cellM24.Value = "2026-02-24"
$cellM24.ClearFormats()

# --- 3) Add a new "Komisyon" sheet (İLÇE / BİRİM / KOMİSYON GÖREVLİLERİ) at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsKomisyon = $wb.Worksheets.Add($null, $lastSheet)
$wsKomisyon.Name = "Komisyon"

$wsKomisyon.Cells.Item(1, 1).Value = "İLÇE"
$wsKomisyon.Cells.Item(1, 2).Value = "BİRİM"
$wsKomisyon.Cells.Item(1, 3).Value = "KOMİSYON GÖREVLİLERİ"

$komisyonData = @(
    @("Akdeniz", "AKDAM"),
    @("Akdeniz", "ESENLİ"),
    @("Akdeniz", "HEBİLLİ"),
    @("Toroslar", "BEKİRALANI"),
    @("Toroslar", "ÇELEBİLİ"),
    @("Toroslar", "DARISEKİSİ"),
    @("Toroslar", "DÜĞDÜÖREN"),
    @("Toroslar", "EVCİLİ"),
    @("Toroslar", "KORUCULAR"),
    @("Toroslar", "MUSALI"),
    @("Tarsus", "GÖÇÜK"),
    @("Tarsus", "KIZILÇUKUR"),
    @("Tarsus", "KARAKÜTÜK"),
    @("Tarsus", "ESKİŞEHİR"),
    @("Tarsus", "KERİMLER"),
    @("Tarsus", "TAŞÇILI"),
    @("Tarsus", "İNCİRGEDİĞİ"),
    @("Tarsus", "KADELLİ"),
    @("Tarsus", "CİNKÖY"),
    @("Tarsus", "İNKÖY"),
    @("Tarsus", "OLUKKOYAĞI"),
    @("Mezitli", "BOZÖN"),
    @("Toroslar", "BELENKEŞLİK"),
    @("Toroslar", "DORUKLU"),
    @("Toroslar", "GÖZNE"),
    @("Toroslar", "KAŞLI"),
    @("Toroslar", "RESULKÖY"),
    @("Yenişehir", "EMİRLER"),
    @("Anamur", "AŞAĞIKÜKÜR"),
    @("Anamur", "BOZDOĞAN"),
    @("Anamur", "BOĞUNTU"),
    @("Anamur", "ÇAMLIPINAR"),
    @("Anamur", "ÇAMLIPINARALANI"),
    @("Anamur", "ÇATALOLUK"),
    @("Anamur", "ÇUKURABANOZ"),
    @("Anamur", "DEMİRÖREN"),
    @("Anamur", "GÜNGÖREN"),
    @("Anamur", "KALINÖREN"),
    @("Anamur", "KARAÇUKUR"),
    @("Anamur", "KARALARBAHŞİŞ"),
    @("Anamur", "KARAAĞA"),
    @("Anamur", "KORUCUK"),
    @("Anamur", "MALAKLAR"),
    @("Anamur", "NASRADDİN"),
    @("Anamur", "ORMANCIK"),
    @("Anamur", "ORTAKÖY"),
    @("Anamur", "SUGÖZÜ"),
    @("Anamur", "SARIAĞAÇ"),
    @("Anamur", "SARIDANA"),
    @("Anamur", "YUKARIKÜKÜR"),
    @("Aydıncık", "HACIBAHATTİN"),
    @("Bozyazı", "BAHÇEKOYAĞI"),
    @("Bozyazı", "BEYRELİ"),
    @("Bozyazı", "ELMAKUZU"),
    @("Bozyazı", "GÖZCE"),
    @("Bozyazı", "GÖZSÜZCE"),
    @("Bozyazı", "KÖMÜRLÜ"),
    @("Bozyazı", "LENGER"),
    @("Bozyazı", "KIZILCA"),
    @("Bozyazı", "TEKEDÜZÜ"),
    @("Bozyazı", "DEREKÖY"),
    @("Gülnar", "BÜYÜKECELİ"),
    @("Gülnar", "ŞEYHÖMER"),
    @("Gülnar", "ZEYNE(SÜTLÜCE)"),
    @("Silifke", "ALTINKUM"),
    @("Silifke", "ARKARASI"),
    @("Silifke", "ATİK"),
    @("Silifke", "BOYNUİNCELİ"),
    @("Silifke", "BURUNUCU"),
    @("Silifke", "ÇELTİKÇİ"),
    @("Silifke", "GÜLÜMPAŞALI"),
    @("Silifke", "HIRMANLI"),
    @("Silifke", "CILBAYIR"),
    @("Silifke", "İMAMUŞAĞI"),
    @("Silifke", "KAVAK"),
    @("Silifke", "KURTULUŞ"),
    @("Silifke", "MARA"),
    @("Silifke", "NASRULLAH"),
    @("Silifke", "SEYDİLİ"),
    @("Silifke", "SÖKÜN"),
    @("Silifke", "ULUGÖZ"),
    @("Silifke", "YEĞENLİ"),
    @("Yenişehir", "Turunçlu"),
    @("Erdemli", "Pınarbaşı"),
    @("Toroslar", "Tırtar"),
    @("Toroslar", "Arslanköy"),
    @("Mut", "Hacınuhlu"),
    @("Silifke", "Kızılisalı"),
    @("Anamur", "Lale"),
    @("Bozyazı", "Derebaşı")
)

$r = 2
foreach ($row in $komisyonData) {
    $wsKomisyon.Cells.Item($r, 1).Value = $row[0]
    $wsKomisyon.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Column C ("KOMİSYON GÖREVLİLERİ") is not yet assigned for any village - keep it as an
# (empty-string) text cell, matching the rest of the sheet, rather than leaving it truly blank.
$wsKomisyon.Range("C2:C91").Formula = '=""'
